$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new Mac-Address rows (31 and 32), following the same pattern
# as the existing rows in the table.
$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 110030
$ws.Cells.Item(31, 3).Value = "eng"
$ws.Cells.Item(31, 4).Value = $true
$ws.Cells.Item(31, 5).Value = "superadmin"
$ws.Cells.Item(31, 6).Value = "now()"
$ws.Cells.Item(31, 7).Value = "now()"

$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 110031
$ws.Cells.Item(32, 3).Value = "eng"
$ws.Cells.Item(32, 4).Value = $true
$ws.Cells.Item(32, 5).Value = "superadmin"
$ws.Cells.Item(32, 6).Value = "now()"
$ws.Cells.Item(32, 7).Value = "now()"

# Update the view so the newly added rows are visible, matching the
# author's final selection/scroll position.
$ws.Range("E28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
